# Updated VisitList Test Plan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Insert two new columns so the existing FIN/Enc* columns shift right into
# their new positions (this carries both content and per-column styling
# along with the shift, matching native Excel column-insert behaviour):
#   - new column F: will hold "FINA" (old F=PatientB -> G, old G=FIN -> H)
#   - new column I (after the relocated FIN column H): will hold
#     "RelationType" (old H..S = Enc1..Enc3 fields -> shift to J..U)
$ws.Columns("F:F").Insert()
$ws.Columns("I:I").Insert()

# --- Row 1 headers: fill in the three new header cells ---
$ws.Range("F1").Value = "FINA"
$ws.Range("H1").Value = "FINB"
$ws.Range("I1").Value = "RelationType"

# --- Row 2 data ---
# Preserve the date-style formatting (originally on "Chest pain", now at
# J2 after the column shift) by copying it onto I2 ("Admitting Physician")
# before J2's own formatting is cleared back to plain/default.
$ws.Range("J2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J2").ClearFormats()

$ws.Range("F2").Value = 1234567
$ws.Range("G2").Value = "Kheang, NoEncntr2"
$ws.Range("H2").Value = 20001379
$ws.Range("I2").Value = "Admitting Physician"
$ws.Range("J2").Value = "Chest pain"

# --- Column widths (per target layout) ---
$ws.Columns("G:H").ColumnWidth = 22.1640625
$ws.Columns("I:I").ColumnWidth = 24.5
$ws.Columns("J:J").ColumnWidth = 15.6640625
$ws.Columns("L:L").ColumnWidth = 12.5
$ws.Columns("N:N").ColumnWidth = 13.5
$ws.Columns("O:P").ColumnWidth = 12.83203125
$ws.Columns("Q:Q").ColumnWidth = 17.6640625
$ws.Columns("R:R").ColumnWidth = 15.33203125
$ws.Columns("T:T").ColumnWidth = 13
$ws.Columns("U:U").ColumnWidth = 18.6640625

# --- Sheet view / selection ---
$ws.Activate()
$ws.Range("A1:A1048576").Select()

# --- Workbook window position/size ---
$excel.Left = 37760
$excel.Top = 5000
$excel.Height = 16300
